$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 13:20"

# Update Navarra row (row 9)
$ws.Range("B9").Value = 2305
$ws.Range("C9").Value = 192
$ws.Range("D9").Value = 2000
$ws.Range("E9").Value = 113

# Update Melilla row (row 58)
$ws.Range("B58").Value = 54
$ws.Range("D58").Value = 53
